$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(" Franco", " Dallara", " IT22G0911802828638328716550")
    ,@(" Massimo", " Ponti", " IT18W8034292104481328055001")
    ,@(" Donato", " Berr�", " IT46C0861777507940288399049")
    ,@(" Irma", " Grossi", " IT52Y0357680245497703784766")
    ,@(" Alessia", " Inzaghi", " IT87S3236428872725803131321")
    ,@(" Fernanda", " Castellitto", " IT77G0817953418013202089394")
    ,@(" Achille", " Balotelli", " IT21A8964054634106625019577")
    ,@(" Giorgio", " Alighieri", " IT14J1519202636883865277705")
    ,@(" Antonina", " Faranda", " IT07Q5630412441332153455987")
    ,@(" Rosario", " Segni", " IT38O7106928503256334884185")
    ,@(" Gastone", " Bernardini", " IT24X2658564212392332035647")
    ,@(" Coluccio", " Viviani", " IT06V7548083700894609986605")
    ,@(" Gabriella", " Varano", " IT03V2159643755478463515322")
    ,@(" Pomponio", " Quasimodo", " IT62Z0418512675125175253843")
    ,@(" Ugolino", " Tomasini", " IT46K5437180682756522779430")
    ,@(" Bianca", " Bongiorno", " IT33Q5910847349340989270366")
    ,@(" Micheletto", " Renzi", " IT92X9428267733656148257548")
    ,@(" Benedetto", " Federici", " IT45J6523854350322959056872")
    ,@(" Domenico", " Barese", " IT14Y5014754707952219213232")
    ,@(" Nanni", " Bragadin", " IT47V9390981684727103408033")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

